# Applies the "23andMe achievements reordering" edit:
#  1. Adds a manual line break (<w:r><w:br/></w:r>) right after the
#     "11/2018 - 08/2022 (Sunnyvale, CA)" run, before the bookmarkEnd
#     that closes the "sunnyvale-ca" bookmark.
#  2. Re-orders (and lightly re-words) the six achievement bullet points
#     that follow, so that:
#        old #4 (RAL)        -> new #1
#        old #5 (MLFlow)     -> new #2 (re-worded)
#        old #6 (IBD graph)  -> new #3
#        old #1 (SNP ETL)    -> new #4
#        old #2 (diabetes)   -> new #5 (re-worded)
#        old #3 (automated)  -> new #6

$d = $word.ActiveDocument

function Find-And-Replace($oldText, $newText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Could not find text: $oldText"
    }
}

# ---------------------------------------------------------------------
# Step 1: insert a line break after the Sunnyvale date/location line,
# keeping it inside the "sunnyvale-ca" bookmark (i.e. before bookmarkEnd).
# ---------------------------------------------------------------------

$findRng = $d.Content
$found = $findRng.Find.Execute("11/2018 - 08/2022 (Sunnyvale, CA)")
if (-not $found) {
    throw "Could not find Sunnyvale heading line"
}

$startPos = $findRng.Start
$findRng.Collapse(0)
$findRng.InsertAfter([char]11)   # vertical-tab char => manual line break (<w:br/>)
$endPos = $findRng.End

# Locate the bookmark that wraps this run so we can re-anchor it across
# the newly inserted break run too (matches the target XML layout where
# bookmarkEnd comes after the <w:br/> run).
$targetBookmark = $null
for ($i = 1; $i -le $d.Bookmarks.Count; $i++) {
    $b = $d.Bookmarks.Item($i)
    if ($b.Name -eq "sunnyvale-ca") {
        $targetBookmark = $b
        break
    }
}
if ($targetBookmark -ne $null) {
    $targetBookmark.Delete()
    $newBookmarkRange = $d.Range($startPos, $endPos)
    $d.Bookmarks.Add("sunnyvale-ca", $newBookmarkRange)
}

# ---------------------------------------------------------------------
# Step 2: reorder / reword the six bullet points.
# Replace each old bullet's text with a unique placeholder token first
# (so that substrings which re-appear elsewhere don't get double-replaced),
# then swap the placeholders for their final text.
# ---------------------------------------------------------------------

$oldTexts = @(
    "Built a large-scale feature engineering ETL pipeline for imputed SNPs (~10 million samples x ~1 million SNPs) using AWS Batch, Metaflow, AWS Glue, and AWS Athena enabling creation of higher quality GWAS and Polygenic Risk Score (PRS) ML models.",
    "Built improved models for type 2 diabetes and Coronary Artery Disease by building model stacking into production PRS pipelines, improving the sensitivity and specificity of 23andMe tests for tens of thousands of customers.",
    "Automated performance metric report generation for all polygenic risk score classifiers leveraging MLFlow artifact storage and headless jupyter execution, reducing researcher time spent on analysis from days to minutes.",
    "Developed and deployed (using MLFlow + AWS Fargate) Recent Ancestor Locations (RAL) - a high precision, high recall country matching algorithm which serves >15 million customers worldwide.",
    "Piloted adoption of MLFlow for experiment tracking and model registry, additionally building completely automated realtime performance metric reporting, eliminating a key source of pipeline fragmentation and redundancy.",
    "Improved graph-based techniques for unsupervised identification of populations by genetically based identity-by-descent (IBD) family relationship, demonstrating an effective way to segment sub-populations (graph community detection) in Mexico and the United Kingdom in an semi-unsupervised manner."
)

$newTexts = @(
    "Created and deployed into production Recent Ancestor Locations (RAL) - a high precision, high recall country matching algorithm which serves >15 million customers worldwide.",
    "Piloted adoption of MLFlow for experiment tracking and model registration, additionally building completely automated realtime performance metric reporting saving researchers hours of time spent on analytics.",
    "Improved graph-based techniques for unsupervised identification of populations by genetically based identity-by-descent (IBD) family relationship, demonstrating an effective way to segment sub-populations (graph community detection) in Mexico and the United Kingdom in an semi-unsupervised manner.",
    "Built a large-scale feature engineering ETL pipeline for imputed SNPs (~10 million samples x ~1 million SNPs) using AWS Batch, Metaflow, AWS Glue, and AWS Athena enabling creation of higher quality GWAS and Polygenic Risk Score (PRS) ML models.",
    "Developed improved models for type 2 diabetes and Coronary Artery Disease by building and evaluating model stacking ensembles into production PRS pipelines, improving the sensitivity and specificity of 23andMe tests for tens of thousands of customers.",
    "Automated performance metric report generation for all polygenic risk score classifiers leveraging MLFlow artifact storage and headless jupyter execution, reducing researcher time spent on analysis from days to minutes."
)

# Pass 1: old text -> unique placeholder token.
for ($i = 0; $i -lt $oldTexts.Length; $i++) {
    $token = "@@REORDER_TOKEN_$i@@"
    Find-And-Replace $oldTexts[$i] $token
}

# Pass 2: placeholder token -> final (new position's) text.
for ($i = 0; $i -lt $newTexts.Length; $i++) {
    $token = "@@REORDER_TOKEN_$i@@"
    Find-And-Replace $token $newTexts[$i]
}

Write-Host "Edit complete"
